# Data update using git
# Apply updated enrollment figures to the "Resumo de Inscrições" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 141

$ws.Range("E10").Value = 602
$ws.Range("F10").Value = 303
$ws.Range("H10").Value = 399

$ws.Range("E11").Value = 386

$ws.Range("F12").Value = 331
$ws.Range("H12").Value = 417

$ws.Range("E16").Value = 211
$ws.Range("F16").Value = 106
$ws.Range("H16").Value = 154

$ws.Range("E17").Value = 107

$ws.Range("F22").Value = 97
$ws.Range("H22").Value = 139

$ws.Range("E25").Value = 291
$ws.Range("F25").Value = 147
$ws.Range("H25").Value = 207

$ws.Range("E27").Value = 345

$ws.Range("E28").Value = 208

$ws.Range("E30").Value = 224

$ws.Range("E32").Value = 192

$ws.Range("E41").Value = 407

$ws.Range("E42").Value = 400

$ws.Range("E45").Value = 158

$ws.Range("E46").Value = 342

$ws.Range("E47").Value = 480

$ws.Range("E48").Value = 230
